$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Zustand" (status) column E for rows that were previously empty
$ws.Range("E15").Value = "In Arbeit"
$ws.Range("E16").Value = "In Arbeit"
$ws.Range("E18").Value = "In Arbeit"

# Update rows that previously said "Bearbeitet"
$ws.Range("E20").Value = "Wartet"
$ws.Range("E21").Value = "In Arbeit"
$ws.Range("E22").Value = "In Arbeit"
$ws.Range("E23").Value = "In Arbeit"

# Update the active selection to match the saved view state
$ws.Range("E24").Select()
